# Atualiza as datas de expiração (coluna B) para refletir a nova janela
# de antecedência de aviso, e move a seleção para B10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value  = 45284
$ws.Range("B3").Value  = 45285
$ws.Range("B4").Value  = 45286
$ws.Range("B5").Value  = 45287
$ws.Range("B6").Value  = 45288
$ws.Range("B7").Value  = 45289
$ws.Range("B8").Value  = 45290
$ws.Range("B9").Value  = 45291
$ws.Range("B10").Value = 45292

$ws.Range("B10").Select()
